$d = $word.ActiveDocument

# Remove the "Tempat Keberangkatan dan Tujuan" bullet paragraph entirely
# (its two runs + paragraph mark), leaving "Data Penerbangan" immediately
# following "Data Pesanan" in the bullet list.
$find1 = $d.Content
$found1 = $find1.Find.Execute("Tempat Keberangkatan dan Tujuan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $d.Range($find1.Start, $find1.End + 1).Delete()
}

# Remove the trailing "Data Penumpang" bullet paragraph that directly
# follows "Data Penerbangan" (the earlier "Data Penumpang" bullet, at the
# top of the same list, must stay untouched).
$find2 = $d.Content
$found2 = $find2.Find.Execute("Data Penerbangan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $afterRange = $d.Range($find2.End, $d.Content.End)
    $found3 = $afterRange.Find.Execute("Data Penumpang", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found3) {
        $d.Range($afterRange.Start, $afterRange.End + 1).Delete()
    }
}
